$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new value is a plain string (non-numeric text, e.g. dotted
# thousands-style prices, or percentage strings) -- a direct .Value assignment
# keeps Excel from reinterpreting them as numbers. ---
$ws.Range("D2").Value = '67.068.60'
$ws.Range("E2").Value = '  -0.78%  '
$ws.Range("D3").Value = '3.512.64'
$ws.Range("E3").Value = '  +0.77%  '
$ws.Range("E4").Value = '  -0.31%  '
$ws.Range("E5").Value = '  +0.84%  '
$ws.Range("E6").Value = '  -0.88%  '
$ws.Range("D7").Value = '3.511.15'
$ws.Range("E7").Value = '  +0.85%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("E9").Value = '  -1.20%  '
$ws.Range("E10").Value = '  -0.57%  '
$ws.Range("E11").Value = '  +7.50%  '
$ws.Range("E12").Value = '  -1.40%  '
$ws.Range("E13").Value = '  +1.10%  '
$ws.Range("D14").Value = '4.106.27'
$ws.Range("E14").Value = '  +0.74%  '
$ws.Range("E15").Value = '  -0.81%  '
$ws.Range("D16").Value = '3.516.15'
$ws.Range("E16").Value = '  +0.87%  '
$ws.Range("D17").Value = '67.123.99'
$ws.Range("E17").Value = '  -0.71%  '
$ws.Range("E18").Value = '  -0.11%  '
$ws.Range("E19").Value = '  +9.03%  '
$ws.Range("E20").Value = '  -1.39%  '
$ws.Range("E21").Value = '  +0.67%  '
$ws.Range("E22").Value = '  -1.64%  '
$ws.Range("E23").Value = '  -2.36%  '
$ws.Range("E24").Value = '  +1.16%  '
$ws.Range("D25").Value = '3.648.91'
$ws.Range("E25").Value = '  +0.50%  '
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("E27").Value = '  -2.43%  '
$ws.Range("E28").Value = '  -1.29%  '
$ws.Range("E29").Value = '  -4.48%  '
$ws.Range("E30").Value = '  +1.20%  '
$ws.Range("E31").Value = '  -4.03%  '
$ws.Range("E32").Value = '  +0.12%  '
$ws.Range("E33").Value = '  -1.00%  '
$ws.Range("E34").Value = '  +0.25%  '
$ws.Range("E35").Value = '  -2.84%  '
$ws.Range("E36").Value = '  -1.75%  '
$ws.Range("E37").Value = '  +0.90%  '
$ws.Range("E38").Value = '  +0.01%  '
$ws.Range("E39").Value = '  -0.29%  '
$ws.Range("E40").Value = '  +0.68%  '
$ws.Range("E41").Value = '  -0.13%  '
$ws.Range("E42").Value = '  +0.41%  '
$ws.Range("E43").Value = '  -10.45%  '
$ws.Range("E44").Value = '  +0.05%  '
$ws.Range("E46").Value = '  -7.00%  '
$ws.Range("E47").Value = '  -3.38%  '
$ws.Range("E48").Value = '  -1.52%  '
$ws.Range("E49").Value = '  -1.00%  '
$ws.Range("E50").Value = '  +0.64%  '
$ws.Range("E51").Value = '  -1.75%  '

# --- Cells whose new value LOOKS like a plain number (e.g. "8.07"). Excel would
# silently store these as numeric cells, losing the original text-cell nature of
# the column. Force text formatting for the assignment, then restore the default
# "Normal" style so no stray style index is left on the cell. ---
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '609.72'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.34'
$ws.Range("D6").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '8.07'
$ws.Range("D11").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '31.56'
$ws.Range("D15").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.81'
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.38'
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.37'
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '436.87'
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '79.72'
$ws.Range("D24").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.77'
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.23'
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.52'
$ws.Range("D30").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '25.58'
$ws.Range("D34").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.03'
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.999'
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '176.11'
$ws.Range("D40").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '46.21'
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '27.98'
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.24'
$ws.Range("D47").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.995'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.247'
$ws.Range("D51").Style = "Normal"
